# The document contains two "<figure>" markup paragraphs. We need the
# second one (the one immediately followed by the paragraph that holds
# "<id>fig_p167r_2</id>"), so walk the paragraphs and disambiguate by
# looking at what follows.
$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "<figure>`r") {
        $nxt = $p.Next()
        if ($nxt -ne $null -and $nxt.Range.Text -eq "<id>fig_p167r_2</id>`r") {
            $target = $p
        }
    }
}

if ($target -eq $null) {
    throw "Could not locate the target '<figure>' paragraph preceding fig_p167r_2"
}

# That paragraph currently holds two runs: the "<figure>" text run, and a
# trailing empty run. Split it into two paragraphs: the first keeps the
# "<figure>" run (and gains paragraph-mark run formatting of Courier New,
# blue, 9pt), the second (newly inserted) takes over the previously
# trailing empty run and carries the same paragraph-mark formatting.
# Doing this via InsertXML (instead of InsertParagraphAfter) lets us set
# the exact paragraph-mark rPr and keep the formerly-trailing empty run
# attached to the new paragraph, matching the target OOXML precisely.

$pPrXml = '<w:widowControl w:val="0"/>' + `
          '<w:pBdr>' + `
            '<w:top w:space="0" w:sz="0" w:val="nil"/>' + `
            '<w:left w:space="0" w:sz="0" w:val="nil"/>' + `
            '<w:bottom w:space="0" w:sz="0" w:val="nil"/>' + `
            '<w:right w:space="0" w:sz="0" w:val="nil"/>' + `
            '<w:between w:space="0" w:sz="0" w:val="nil"/>' + `
          '</w:pBdr>' + `
          '<w:shd w:fill="auto" w:val="clear"/>' + `
          '<w:contextualSpacing w:val="0"/>' + `
          '<w:rPr>' + `
            '<w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/>' + `
            '<w:color w:val="0000ff"/>' + `
            '<w:sz w:val="18"/>' + `
            '<w:szCs w:val="18"/>' + `
          '</w:rPr>'

$firstPara = '<w:p w:rsidR="00000000" w:rsidDel="00000000" w:rsidP="00000000" w:rsidRDefault="00000000" w:rsidRPr="00000000" w14:paraId="00000015">' + `
               '<w:pPr>' + $pPrXml + '</w:pPr>' + `
               '<w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000">' + `
                 '<w:rPr>' + `
                   '<w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/>' + `
                   '<w:color w:val="0000ff"/>' + `
                   '<w:sz w:val="18"/>' + `
                   '<w:szCs w:val="18"/>' + `
                   '<w:rtl w:val="0"/>' + `
                 '</w:rPr>' + `
                 '<w:t xml:space="preserve">&lt;figure&gt;</w:t>' + `
               '</w:r>' + `
             '</w:p>'

$secondPara = '<w:p>' + `
                '<w:pPr>' + $pPrXml + '</w:pPr>' + `
                '<w:r>' + `
                  '<w:rPr>' + `
                    '<w:rtl w:val="0"/>' + `
                  '</w:rPr>' + `
                '</w:r>' + `
              '</w:p>'

$bodyXml = $firstPara + $secondPara

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
              '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
                '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
                  '<pkg:xmlData>' + `
                    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
                      '<w:body>' + $bodyXml + '</w:body>' + `
                    '</w:document>' + `
                  '</pkg:xmlData>' + `
                '</pkg:part>' + `
              '</pkg:package>'

$target.Range.InsertXML($packageXml)
